$d = $word.ActiveDocument

# Add the June 7th diary sentence to the end of the last paragraph's run.
[void]$d.Content.Find.Execute("晴，今天是高考的一天，上午考语文，下午考英语。", $true, $false, $false, $false, $false,
                         $true, 1, $false, "晴，今天是高考的一天，上午考语文，下午考英语。今天天气真不错。", 2)

# Mark the "Default Paragraph Font" and "Normal Table" built-in styles as
# recommended Quick Styles (w:qFormat), matching the style-panel refresh
# Word performs whenever it resaves this style sheet.
[void]($d.Styles("Default Paragraph Font").QuickStyle = $true)
[void]($d.Styles("Normal Table").QuickStyle = $true)
